$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.98"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.875.93"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.86"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.03"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08049"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.12"
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").Value = "1.881.22"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.961"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.111"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.05"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06692"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001043"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "27.669.42"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.524"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.309"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "2.110.84"
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.54"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.099"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.555"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.75"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9801"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09473"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.444"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.602"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.350"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06115"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.116"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5978"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1894"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.260"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5708"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.18"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.396"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.942"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06919"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.51"
$ws.Range("E50").Value = "  +6.06%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.070"
$ws.Range("E51").Value = "  +1.84%  "
